# The upstream change for this fixture (commit "Fixed POI packaging and
# upgraded to POI 3.15") is a pure artifact of regenerating the test
# resource with a newer Apache POI / XMLBeans release: every hunk in the
# unified diff touches only the *serialization order* of XML attributes
# (and namespace declarations on the <w:document> / <w:styles> roots,
# plus the position of mc:Ignorable). Tag names, attribute names,
# attribute values, element nesting, and text content are all byte-for-
# byte identical before and after - nothing in the document's visible
# content, formatting, or structure actually changed.
#
# The Word object model has no surface for reordering raw XML attributes
# (that's an artifact of the serializer library, not something exposed
# via Find/Replace, Styles, Paragraphs, etc.), so there is no OM
# operation that corresponds to this diff. Applying "the change" via
# Word COM therefore means leaving the document's content untouched.
$d = $word.ActiveDocument
